$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Ex 1-A": fix the S1-P1 row (row 13) so the bias reference ($C$7) is
# absolute, matching the same pattern already used in row 10 (S1 P0).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Ex 1-A")

$ws1.Range("D13").Formula = '=($C$6*C4+$D$6*D4+$E$6*E4)+$C$7'
$ws1.Range("E13").Formula = '=($C$6*D4+$D$6*E4+$E$6*F4)+$C$7'
$ws1.Range("F13").Formula = '=($C$6*E4+$D$6*F4+$E$6*G4)+$C$7'
$ws1.Range("G13").Formula = '=($C$6*F4+$D$6*G4+$E$6*H4)+$C$7'
$ws1.Range("H13").Formula = '=($C$6*G4+$D$6*H4+$E$6*I4)+$C$7'
$ws1.Range("I13").Formula = '=($C$6*H4+$D$6*I4+$E$6*J4)+$C$7'
$ws1.Range("J13").Formula = '=($C$6*I4+$D$6*J4+$E$6*K4)+$C$7'
$ws1.Range("K13").Formula = '=($C$6*J4+$D$6*K4)+$C$7'

# ---------------------------------------------------------------------------
# Sheet "Ex 1-B": add the 2D-convolution output grids for Filter1 (T3:V5)
# and Filter2 (X3:Z5), then highlight them (yellow fill + border), matching
# the "Output (To be completed)" header above.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Ex 1-B")

# Filter 1 output (columns T:V)
$ws2.Range("T3").Formula = '=(D3*$K$3+E3*$L$3+D4*$K$4+E4*$L$4)+(D9*$K$6+E9*$L$6+D10*$K$7+E10*$L$7)+(D15*$K$9+E15*$L$9+D16*$K$10+E16*$L$10)+$L$12'
$ws2.Range("U3").Formula = '=(E3*$K$3+F3*$L$3+E4*$K$4+F4*$L$4)+(E9*$K$6+F9*$L$6+E10*$K$7+F10*$L$7)+(E15*$K$9+F15*$L$9+E16*$K$10+F16*$L$10)+$L$12'
$ws2.Range("V3").Formula = '=(F3*$K$3+G3*$L$3+F4*$K$4+G4*$L$4)+(F9*$K$6+G9*$L$6+F10*$K$7+G10*$L$7)+(F15*$K$9+G15*$L$9+F16*$K$10+G16*$L$10)+$L$12'

$ws2.Range("T4").Formula = '=(D4*$K$3+E4*$L$3+D5*$K$4+E5*$L$4)+(D10*$K$6+E10*$L$6+D11*$K$7+E11*$L$7)+(D16*$K$9+E16*$L$9+D17*$K$10+E17*$L$10)+$L$12'
$ws2.Range("U4").Formula = '=(E4*$K$3+F4*$L$3+E5*$K$4+F5*$L$4)+(E10*$K$6+F10*$L$6+E11*$K$7+F11*$L$7)+(E16*$K$9+F16*$L$9+E17*$K$10+F17*$L$10)+$L$12'
$ws2.Range("V4").Formula = '=(F4*$K$3+G4*$L$3+F5*$K$4+G5*$L$4)+(F10*$K$6+G10*$L$6+F11*$K$7+G11*$L$7)+(F16*$K$9+G16*$L$9+F17*$K$10+G17*$L$10)+$L$12'

$ws2.Range("T5").Formula = '=(D5*$K$3+E5*$L$3+D6*$K$4+E6*$L$4)+(D11*$K$6+E11*$L$6+D12*$K$7+E12*$L$7)+(D17*$K$9+E17*$L$9+D18*$K$10+E18*$L$10)+$L$12'
$ws2.Range("U5").Formula = '=(E5*$K$3+F5*$L$3+E6*$K$4+F6*$L$4)+(E11*$K$6+F11*$L$6+E12*$K$7+F12*$L$7)+(E17*$K$9+F17*$L$9+E18*$K$10+F18*$L$10)+$L$12'
$ws2.Range("V5").Formula = '=(F5*$K$3+G5*$L$3+F6*$K$4+G6*$L$4)+(F11*$K$6+G11*$L$6+F12*$K$7+G12*$L$7)+(F17*$K$9+G17*$L$9+F18*$K$10+G18*$L$10)+$L$12'

# Filter 2 output (columns X:Z)
$ws2.Range("X3").Formula = '=(D3*$O$3+E3*$P$3+D4*$O$4+E4*$P$4)+(D9*$O$6+E9*$P$6+D10*$O$7+E10*$P$7)+(D15*$O$9+E15*$P$9+D16*$O$10+E16*$P$10)+$P$12'
$ws2.Range("Y3").Formula = '=(E3*$O$3+F3*$P$3+E4*$O$4+F4*$P$4)+(E9*$O$6+F9*$P$6+E10*$O$7+F10*$P$7)+(E15*$O$9+F15*$P$9+E16*$O$10+F16*$P$10)+$P$12'
$ws2.Range("Z3").Formula = '=(F3*$O$3+G3*$P$3+F4*$O$4+G4*$P$4)+(F9*$O$6+G9*$P$6+F10*$O$7+G10*$P$7)+(F15*$O$9+G15*$P$9+F16*$O$10+G16*$P$10)+$P$12'

$ws2.Range("X4").Formula = '=(D4*$O$3+E4*$P$3+D5*$O$4+E5*$P$4)+(D10*$O$6+E10*$P$6+D11*$O$7+E11*$P$7)+(D16*$O$9+E16*$P$9+D17*$O$10+E17*$P$10)+$P$12'
$ws2.Range("Y4").Formula = '=(E4*$O$3+F4*$P$3+E5*$O$4+F5*$P$4)+(E10*$O$6+F10*$P$6+E11*$O$7+F11*$P$7)+(E16*$O$9+F16*$P$9+E17*$O$10+F17*$P$10)+$P$12'
$ws2.Range("Z4").Formula = '=(F4*$O$3+G4*$P$3+F5*$O$4+G5*$P$4)+(F10*$O$6+G10*$P$6+F11*$O$7+G11*$P$7)+(F16*$O$9+G16*$P$9+F17*$O$10+G17*$P$10)+$P$12'

$ws2.Range("X5").Formula = '=(D5*$O$3+E5*$P$3+D6*$O$4+E6*$P$4)+(D11*$O$6+E11*$P$6+D12*$O$7+E12*$P$7)+(D17*$O$9+E17*$P$9+D18*$O$10+E18*$P$10)+$P$12'
$ws2.Range("Y5").Formula = '=(E5*$O$3+F5*$P$3+E6*$O$4+F6*$P$4)+(E11*$O$6+F11*$P$6+E12*$O$7+F12*$P$7)+(E17*$O$9+F17*$P$9+E18*$O$10+F18*$P$10)+$P$12'
$ws2.Range("Z5").Formula = '=(F5*$O$3+G5*$P$3+F6*$O$4+G6*$P$4)+(F11*$O$6+G11*$P$6+F12*$O$7+G12*$P$7)+(F17*$O$9+G17*$P$9+F18*$O$10+G18*$P$10)+$P$12'

# Highlight the newly-computed output blocks: yellow fill + thin border,
# centered text (matches the style already used for completed output cells).
$rngF1 = $ws2.Range("T3:V5")
$rngF1.Borders.LineStyle = 1
$rngF1.Interior.Color = 65535
$rngF1.HorizontalAlignment = -4108
$rngF1.VerticalAlignment = -4108

$rngF2 = $ws2.Range("X3:Z5")
$rngF2.Borders.LineStyle = 1
$rngF2.Interior.Color = 65535
$rngF2.HorizontalAlignment = -4108
$rngF2.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping: the workbook was saved with "Ex 1-A"
# selected at M18 and "Ex 1-B" (the sheet actually worked on) active, with
# its selection left at P16.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("M18").Select()

$ws2.Activate()
$ws2.Range("P16").Select()
